$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in new cell text content (columns A-D) for rows 2-15 ---
$ws.Range("A2").Value = '(SEO ou accessiblité ?)'
$ws.Range("B2").Value = 'Format images'
$ws.Range("C2").Value = 'Les formats d''image comme JPEG 2000, JPEG XR et WebP proposent souvent une meilleure compression que les formats PNG ou JPEG. Par conséquent, les téléchargements sont plus rapides et la consommation de données est réduite.'
$ws.Range("D2").Value = 'compresser les images'
$ws.Range("F2").Value = 'https://web.dev/uses-webp-images'

$ws.Range("B3").Value = 'Images lourde'
$ws.Range("C3").Value = 'Les images optimisées se chargent plus rapidement et consomment moins de données mobiles.'
$ws.Range("D3").Value = 'Les images optimisées se chargent plus rapidement et consomment moins de données mobiles.'
$ws.Range("F3").Value = 'https://web.dev/uses-optimized-images'

$ws.Range("B4").Value = 'taille des images'
$ws.Range("C4").Value = 'Diffusez des images de taille appropriée afin d''économiser des données mobiles et de réduire le temps de chargement.'
$ws.Range("D4").Value = 'Diffusez des images de taille appropriée afin d''économiser des données mobiles et de réduire le temps de chargement.'
$ws.Range("F4").Value = 'https://web.dev/uses-responsive-images'

$ws.Range("B5").Value = 'cache'
$ws.Range("C5").Value = 'Une longue durée de vie du cache peut accélérer les visites répétées sur votre page.'
$ws.Range("D5").Value = 'Une longue durée de vie du cache peut accélérer les visites répétées sur votre page.'
$ws.Range("F5").Value = 'https://web.dev/uses-long-cache-ttl'

$ws.Range("B6").Value = 'texte non visible pendant le chargement'
$ws.Range("C6").Value = 'Utilisez la fonction d''affichage de la police CSS afin que le texte soit visible par l''utilisateur pendant le chargement des polices Web.'
$ws.Range("D6").Value = 'Utilisez la fonction d''affichage de la police CSS afin que le texte soit visible par l''utilisateur pendant le chargement des polices Web.'

$ws.Range("B7").Value = 'pas de width ni de height aux images'
$ws.Range("C7").Value = 'Indiquez une largeur et une hauteur explicites sur les éléments d''image afin de réduire les décalages de mise en page et d''améliorer le CLS.'
$ws.Range("D7").Value = 'Indiquez une largeur et une hauteur explicites sur les éléments d''image afin de réduire les décalages de mise en page et d''améliorer le CLS.'

$ws.Range("B8").Value = 'Évitez de créer des chaînes de requêtes critiques'
$ws.Range("C8").Value = 'Les chaînes de demandes critiques ci-dessous vous montrent quelles ressources sont chargées avec une priorité élevée. Envisagez de réduire la longueur des chaînes et la taille de téléchargement des ressources ou de reporter le téléchargement de ressources inutiles afin d''améliorer le chargement des pages.'
$ws.Range("D8").Value = 'Les chaînes de demandes critiques ci-dessous vous montrent quelles ressources sont chargées avec une priorité élevée. Envisagez de réduire la longueur des chaînes et la taille de téléchargement des ressources ou de reporter le téléchargement de ressources inutiles afin d''améliorer le chargement des pages.'

$ws.Range("B9").Value = 'Réduisez au maximum le nombre de requêtes et la taille des transferts'
$ws.Range("C9").Value = 'Pour définir des budgets liés à la quantité et à la taille des ressources de pages, ajoutez un fichier budget.json.'
$ws.Range("D9").Value = 'Pour définir des budgets liés à la quantité et à la taille des ressources de pages, ajoutez un fichier budget.json.'

$ws.Range("B10").Value = 'Élément identifié comme "Largest Contentful Paint"'
$ws.Range("C10").Value = 'Il s''agit de l''élément identifié comme "Largest Contentful Paint" dans la fenêtre d''affichage.'
$ws.Range("D10").Value = 'Il s''agit de l''élément identifié comme "Largest Contentful Paint" dans la fenêtre d''affichage.'

$ws.Range("B11").Value = 'Éviter les changements de mise en page importants'
$ws.Range("C11").Value = 'Ces éléments DOM contribuent en grande partie au CLS de la page.'
$ws.Range("D11").Value = 'Ces éléments DOM contribuent en grande partie au CLS de la page.'

$ws.Range("B12").Value = 'fichier non utiliser'
$ws.Range("C12").Value = 'il y a un fichier javascript et css non utilisé lors du chargement'

$ws.Range("A13").Value = 'SEO'
$ws.Range("B13").Value = 'pas de meta description'
$ws.Range("F13").Value = 'https://web.dev/meta-description/'

$ws.Range("B14").Value = 'Évitez d''énormes charges utiles de réseau'
$ws.Range("C14").Value = 'Les charges utiles des grands réseaux coûtent de l''argent réel aux utilisateurs et sont fortement corrélées aux délais de chargement interminables'
$ws.Range("D14").Value = 'Les charges utiles des grands réseaux coûtent de l''argent réel aux utilisateurs et sont fortement corrélées aux délais de chargement interminables'
$ws.Range("F14").Value = 'https://web.dev/total-byte-weight'

$ws.Range("A15").Value = 'accessibilité'
$ws.Range("B15").Value = 'html n''a pas de lang'
$ws.Range("F15").Value = 'https://web.dev/html-lang-valid/'

# --- Add hyperlinks for column F, rows 10 down to 6 (so relationship ids come out in this order) ---
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://web.dev/lighthouse-largest-contentful-paint') | Out-Null
$ws.Range("F10").Value = 'https://web.dev/lighthouse-largest-contentful-paint'

$ws.Hyperlinks.Add($ws.Range("F9"), 'https://web.dev/use-lighthouse-for-performance-budgets') | Out-Null
$ws.Range("F9").Value = 'https://web.dev/use-lighthouse-for-performance-budgets'

$ws.Hyperlinks.Add($ws.Range("F8"), 'https://web.dev/critical-request-chains') | Out-Null
$ws.Range("F8").Value = 'https://web.dev/critical-request-chains'

$ws.Hyperlinks.Add($ws.Range("F7"), 'https://web.dev/optimize-cls/?utm_source=lighthouse&utm_medium=unknown', 'images-without-dimensions', $null, 'https://web.dev/optimize-cls/?utm_source=lighthouse&utm_medium=unknown#images-without-dimensions') | Out-Null
$ws.Range("F7").Value = 'https://web.dev/optimize-cls#images-without-dimensions'

$ws.Hyperlinks.Add($ws.Range("F6"), 'https://web.dev/font-display') | Out-Null
$ws.Range("F6").Value = 'https://web.dev/font-display'

# --- Set the active selection to match the authored workbook (C12) ---
$ws.Range("C12").Select() | Out-Null

